# otras mejoreas al modelo Balancear entre números de alta confianza y diversificación
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Add 4 new columns (FZ, GA, GB, GC) of data for rows 2-7
$data = @{
    2 = @(3, 1, 0, 8)
    3 = @(13, 7, 5, 11)
    4 = @(15, 11, 6, 12)
    5 = @(19, 20, 12, 14)
    6 = @(32, 22, 20, 15)
    7 = @(36, 31, 30, 29)
}

$cols = @("FZ", "GA", "GB", "GC")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $row
        $ws.Range($addr).Value = $vals[$i]
    }
}

# Update selection / active cell to match the new last column
$ws.Range("GC2:GC7").Select()
